# Remove the post entry in row 761 ("幸運にも母親とその子２頭を５日間観察できた…")
# and shift every row below it up by one (Excel renumbers A/B/C refs and
# shrinks the sheet's used range automatically).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(761).Delete()
